# Add two new rows (21 and 22) to the bonification tracking sheet,
# continuing the existing data table with the same layout/formatting
# as the preceding rows (copy row 20's formatting down, then fill in
# the new values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: Ordem de Servico 613, Bonificacao 5, Tecnico Bruno ---
$ws.Range("A20:E20").Copy()
$ws.Range("A21").Insert("xlShiftDown")

$ws.Range("A21").Value = 19
$ws.Range("B21").Value2 = 45828
$ws.Range("C21").Value = 613
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = "Bruno"

# --- Row 22: Ordem de Servico 613, Bonificacao 5, Tecnico Ronaldo ---
$ws.Range("A20:E20").Copy()
$ws.Range("A22").Insert("xlShiftDown")

$ws.Range("A22").Value = 20
$ws.Range("B22").Value2 = 45828
$ws.Range("C22").Value = 613
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = "Ronaldo"
